$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (prices / 1h volume %) + Litecoin <-> WrappedEther row swap (rows 14-15)

$ws.Range("D2").Value = "29.897.92"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.888.74"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7631"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3123"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.61"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07163"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08511"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.46%  "
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.363"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.18%  "
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.784.42"
$ws.Range("E15").Value = "  -5.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.140"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").Value = "29.746.18"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007794"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9994"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "2.109.42"
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.968"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1614"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.406"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.033"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.472"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.534"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.488"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.095"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05441"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.90%  "
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7421"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9984"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.698"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01945"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.781"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4459"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").Value = "1.104.94"
$ws.Range("E42").Value = "  -4.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.059"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8529"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("E48").Value = "  -2.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.621"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.997"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.22%  "
$ws.Range("D51").Value = "2.003.11"
$ws.Range("E51").Value = "  -1.88%  "
